$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BWP Bootstrap test case rows appended to the EmulatorData sheet (rows 12-17)
$newRows = @(
    @{ Row = 12; A = "DualCFCeiling";      B = 11; C = 952; D = 1.5; E = 80; F = 2 },
    @{ Row = 13; A = "DualCFFlat";         B = 12; C = 952; D = 1.5; E = 10; F = 2 },
    @{ Row = 14; A = "DualCFPercentage";   B = 13; C = 952; D = 1.5; E = 60; F = 2 },
    @{ Row = 15; A = "SingleCFCeiling";    B = 14; C = 951; D = 1.5; E = 80; F = 2 },
    @{ Row = 16; A = "SingleCFFlat";       B = 15; C = 951; D = 1.5; E = 10; F = 2 },
    @{ Row = 17; A = "SingleCFPercentage"; B = 16; C = 951; D = 1.5; E = 60; F = 2 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = "udf data 1"
    $ws.Range("H$row").Value = "udf data 2"
    $ws.Range("I$row").Value = "udf data 3"
    $ws.Range("J$row").Value = "Sweet"
    $ws.Range("K$row").Value = "Sour"
    $ws.Range("L$row").Value = "udf data 6"
    $ws.Range("M$row").Value = "udf data 7"
    $ws.Range("N$row").Value = "udf data 8"
    $ws.Range("O$row").Value = "udf data 9"
    $ws.Range("P$row").Value = "udf data 10"
}

$ws.Range("A17").Select()
